$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 (which currently holds the "SQL7" label in column B),
# which shifts rows 13-19 up by one.
$ws.Rows.Item(12).Delete()

# After deletion, the former row 19 (数据导出) is now row 18.
# Insert a brand-new row 19 describing a hardware comparison test.
$ws.Cells.Item(19, 1).Value = 6
$ws.Cells.Item(19, 2).Value = "硬件对比"
$ws.Cells.Item(19, 3).Value = "测试环境为SSD硬盘`n生产环境为HDD硬盘`n经测试在测试环境，查询速度明显比生产环境快。`n例如：同样的SQL语句，返回70万条记录在测试环境仅需4分钟，而在生产环境则为9:46；又如，同样的SQL，在测试环境仅需3秒，而在生产环境需要11秒。`n建议：升级生产环境为SSD"

# Column A/B already pick up the correct column-level styles (center/top and
# top alignment respectively) automatically. Column C needs wrap text like
# the other long-text cells in column C.
$ws.Cells.Item(19, 3).WrapText = $true

$ws.Rows.Item(19).RowHeight = 81

# Row 11 ("SQL6") is then relabeled to the combined "SQL6、SQL7".
$ws.Range("B11").Value = "SQL6、SQL7"

# Update selection to match the committed state
$ws.Range("A12:XFD12").Select()
